# Update the "2024" worksheet: a new September transaction arrived, so the
# existing September_Details/September_Date (R:S) entries all move down one
# row and the newest entry is recorded at the top (row 40, right after the
# last already-categorised row). This naturally pushes the trailing
# "Broadband" label (row 158, column A) down into a brand-new row 159 too,
# since it lives below the insertion point.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2024")

# Insert a cell above row 40 in columns R:S only, shifting R40:S<end> (and
# everything below it, including the lone A158 "Broadband" label) down by
# one row.
$ws.Range("R40:S40").Insert()

# Record the newest entry in the freshly opened row.
$ws.Range("R40").Value2 = "balance your axis"
$ws.Range("S40").Value2 = "2024-09-17 13:07:16"

# Make sure the "Broadband" label ends up on row 159 (it is carried there
# automatically by the insert above, but set it explicitly for safety).
$ws.Range("A158").ClearContents()
$ws.Range("A159").Value2 = "Broadband"
